$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PLAYER_DASH_PAUSE_FRAME (row 9) and PLAYER_HP (row 10) values
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 20

# Move the active selection to B11 (matches the new sheetView selection)
$ws.Range("B11").Select()
